# Adds more "dummy" readings data to Sheet1 (rows 46-60), plus a block of
# styled-but-empty placeholder rows (61-69), matching the author's commit
# "added more dummy data".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New glucose-reading / food-log rows that extend the existing table.
# Columns: Before_CBG_Reading, Before_CBG_Measurement, Before_CBG_Uploaded_At,
#          Food_Name, Food_Calorie, Food_Carb, Food_Sugar, Food_Fibre,
#          Food_Uploaded_At, After_CBG_Reading, After_CBG_Measurement,
#          After_CBG_Uploaded_At
$newRows = @(
    @(46, 6.3, "2023-04-03 10:00", "Whole Grain Bread", 80,    20,   4,   3,    "2023-04-03 10:30", 6.4, "2023-04-03 11:15"),
    @(47, 6.4, "2023-04-03 18:00", "Burger",             266,  30.3, 5.2, 1.1,  "2023-04-03 18:30", 6.6, "2023-04-03 19:15"),
    @(48, 6.5, "2023-04-04 10:00", "Nugget",             284.2,11.5, 0,   2,    "2023-04-04 10:30", 6.7, "2023-04-04 11:15"),
    @(49, 6.7, "2023-04-04 18:00", "Burger",             266,  30.3, 5.2, 1.1,  "2023-04-04 18:30", 7,   "2023-04-04 19:15"),
    @(50, 6.8, "2023-04-05 10:00", "Nugget",             284.2,11.5, 0,   2,    "2023-04-05 10:30", 7,   "2023-04-05 11:15"),
    @(51, 6.9, "2023-04-05 18:00", "Burger",             266,  30.3, 5.2, 1.1,  "2023-04-05 18:30", 7.2, "2023-04-05 19:15"),
    @(52, 6.7, "2023-04-06 10:00", "Nugget",             284.2,11.5, 0,   2,    "2023-04-06 10:30", 6.9, "2023-04-06 11:15"),
    @(53, 6.8, "2023-04-06 18:00", "Brown Rice",         248,  51.7, 0.3, 3.23, "2023-04-06 18:30", 6.8, "2023-04-06 19:15"),
    @(54, 6.7, "2023-04-07 10:00", "Nugget",             284.2,11.5, 0,   2,    "2023-04-07 10:30", 6.9, "2023-04-07 11:15"),
    @(55, 6.9, "2023-04-07 18:00", "Brown Rice",         248,  51.7, 0.3, 3.23, "2023-04-07 18:30", 7,   "2023-04-07 19:15"),
    @(56, 6,   "2023-04-08 10:00", "Whole Grain Bread", 80,    20,   4,   3,    "2023-04-08 10:30", 6.1, "2023-04-08 11:15"),
    @(57, 6,   "2023-04-08 18:00", "Burger",             266,  30.3, 5.2, 1.1,  "2023-04-08 18:30", 6.3, "2023-04-08 19:15"),
    @(58, 6.2, "2023-04-09 10:00", "Nugget",             284.2,11.5, 0,   2,    "2023-04-09 10:30", 6.5, "2023-04-09 11:15"),
    @(59, 6.4, "2023-04-09 18:00", "Burger",             266,  30.3, 5.2, 1.1,  "2023-04-09 18:30", 6.7, "2023-04-09 19:15"),
    @(60, 6.4, "2023-04-10 10:00", "Whole Grain Bread", 80,    20,   4,   3,    "2023-04-10 10:30", 6.4, "2023-04-10 11:15")
)

foreach ($row in $newRows) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = "mmolL"
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = "mmolL"
    $ws.Cells.Item($r, 12).Value = $row[10]

    # Date/time stamp columns are kept as plain text (same as the rest of
    # the sheet) instead of being auto-converted to date serials.
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 9).NumberFormat = "@"
    $ws.Cells.Item($r, 12).NumberFormat = "@"

    # All twelve cells in these freshly-entered rows pick up an explicit
    # black font colour (rather than the theme colour used by the older
    # rows), which is how the new rows end up on a distinct style index.
    for ($c = 1; $c -le 12; $c++) {
        $ws.Cells.Item($r, $c).Font.Color = 0
    }
}

# A further block of blank rows (61-69) was added underneath, already
# carrying the same "new" styling (black font / text format) even though
# no values were entered into them yet.
for ($r = 61; $r -le 69; $r++) {
    for ($c = 1; $c -le 12; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Font.Color = 0
        if ($c -eq 3 -or $c -eq 9 -or $c -eq 12) {
            $cell.NumberFormat = "@"
        }
    }
}

$null = $ws.Columns.AutoFit()

$null = $ws.Range("L61").Select()
